$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.895.03'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '2.682.56'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.17'
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.44'
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("E9").Value = '  -3.12%  '
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.366'
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("D13").Value = '3.157.49'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.24'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '62.806.86'
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("D17").Value = '2.683.77'
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.82'
$ws.Range("E18").Value = '  -5.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.59'
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.64'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.21'
$ws.Range("E21").Value = '  -5.26%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.33'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.15'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.41'
$ws.Range("E28").Value = '  +5.55%  '
$ws.Range("D29").Value = '0.0₃0851'
$ws.Range("E29").Value = '  -5.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.72'
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.86'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.46'
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.44'
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '340.08'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.17'
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.928'
$ws.Range("E40").Value = '  -4.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -2.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.33'
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.87'
$ws.Range("E43").Value = '  -3.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.19'
$ws.Range("E44").Value = '  -3.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.616'
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0556'
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.02'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.39'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0969'
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0241'
$ws.Range("E51").Value = '  -3.29%  '
